{"js": "const body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Table rows are 0-indexed; column is always 0 (single-column table).\nconst updates = [\n  [0, \"0M\"],\n  [1, \"0M\"],\n  [2, \"0M\"],\n  [3, \"2339\"],\n  [5, \"0.08911\"],\n  [6, \"0.00846\"],\n  [7, \"0.00297\"],\n  [9, \"0.05209\"],\n  [10, \"0.07973\"],\n  [11, \"3.03292\"],\n  [43, \"99.61\"],\n  [44, \"3.03\"],\n  [45, \"772\"],\n];\n\nfor (const [rowIndex, text] of updates) {\n  const cell = table.getCell(rowIndex, 0);\n  cell.value = text;\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfunction Set-CellText($rowIndex, $text) {\n    $cell = $t.Rows.Item($rowIndex).Cells.Item(1)\n    $cell.Range.Text = $text\n}\n\n# Rows are 1-based in the Word object model.\nSet-CellText 1 \"0M\"\nSet-CellText 2 \"0M\"\nSet-CellText 3 \"0M\"\nSet-CellText 4 \"2339\"\nSet-CellText 6 \"0.08911\"\nSet-CellText 7 \"0.00846\"\nSet-CellText 8 \"0.00297\"\nSet-CellText 10 \"0.05209\"\nSet-CellText 11 \"0.07973\"\nSet-CellText 12 \"3.03292\"\nSet-CellText 44 \"99.61\"\nSet-CellText 45 \"3.03\"\nSet-CellText 46 \"772\"\n"}
